# Budget Plan (version 1).xlsx - update
# Summary of changes:
#   1. Insert a new line item row at row 7:
#        "RC Servo BEC UBEC 3A 5V ( Receiver servo Power supply)" | Qty 1 | $7.99 | =B7*C7
#      (this pushes all the existing rows 7-18 down to 8-19)
#   2. Quantity change: "Cytron 13A, 5-30V Single DC Motor Controller" 3 -> 4
#   3. Quantity change: "Raspberry Pi Zero W" 1 -> 2
#   4. Quantity change: "5" wheels/hubs" 4 -> 6 (this is row 18 after the insert)
#   5. Selection moved to D7

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert the new row for the RC servo BEC part ------------------------
$ws.Rows(7).Insert()

$ws.Range("A7").Value = "RC Servo BEC UBEC 3A 5V ( Receiver servo Power supply)"
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = 7.99
$ws.Range("D7").Formula = "=B7*C7"
$ws.Rows(7).RowHeight = 25.5

# --- 2 & 3. Quantity updates on existing rows --------------------------------
$ws.Range("B3").Value = 4
$ws.Range("B4").Value = 2

# --- 4. Quantity update on the (now shifted) "5 inch wheels/hubs" row -------
$ws.Range("B18").Value = 6

# --- 5. Selection / view -----------------------------------------------------
[void]$ws.Range("D7").Select()
